$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: "В подложке должна располагаться ..." -> "В этом разделе должна
# располагаться ..." split into three runs: "В " / "этом разделе" / " должна
# располагаться информация о правообладателях, ссылка на "
# ---------------------------------------------------------------------------
$find1 = $d.Content
[void]$find1.Find.Execute("В подложке должна располагаться информация", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s1 = $find1.Start

# "В " occupies the first 2 characters, "подложке" the next 8 characters.
$wordRng = $d.Range($s1 + 2, $s1 + 10)
$wordRng.Text = "этом разделе"
$wordRng.Font.Size = 15
$wordRng.Font.Size = 14

# ---------------------------------------------------------------------------
# Part 2: remove "веб-каталога" (and its proofErr wrapper) and rework the
# surrounding text, then split the trailing text into separate runs so that
# "да" + "н" + "н" + "ых" become individual runs.
# ---------------------------------------------------------------------------
$find2 = $d.Content
[void]$find2.Find.Execute(" Разработчика и логотип Разработчика. В случае регистрации сайта в веб-каталога в подложке могут располагаться счетчики данных каталогов.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s2 = $find2.Start
$e2 = $find2.End
$fullRng = $d.Range($s2, $e2)
$fullRng.Text = " Разработчика и логотип Разработчика. В случае регистрации сайта в подложке могут располагаться счетчики данных каталогов."

# Re-locate the tail of the sentence to compute split boundaries.
$find3 = $d.Content
[void]$find3.Find.Execute("могут располагаться счетчики данных каталогов.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailStart = $find3.Start
$tailEnd = $find3.End

$lenD = 31  # length of "могут располагаться счетчики да"
$posD_start = $tailStart
$posD_end = $tailStart + $lenD
$posE_end = $posD_end + 1
$posF_end = $posE_end + 1

$rD = $d.Range($posD_start, $posD_end)
$rD.Font.Size = 15
$rD.Font.Size = 14

$rE = $d.Range($posD_end, $posE_end)
$rE.Font.Size = 15
$rE.Font.Size = 14

$rF = $d.Range($posE_end, $posF_end)
$rF.Font.Size = 15
$rF.Font.Size = 14
